# Automatic update of files.
#
# The underlying records for rows 4-10 rotate down by one: the data that
# used to live in row N now lives in row N+1, and the data that used to
# live in row 10 wraps around into row 4.
#
# Only columns A, B, E, F, G, H, Q, R, AC actually differ row-to-row (every
# other column in this block happens to already hold the same value across
# rows 4-10), so those are the only cells that need to be touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (becomes what row 10 used to hold: "Granticka") ---------------
$ws.Range("A4").Value = 110060427
$ws.Range("B4").Value = 89410
$ws.Range("E4").Value = 5432
$ws.Range("F4").Value = "Granticka"
$ws.Range("G4").Value = "Porodaedalea chrysoloma"
$ws.Range("H4").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q4").Value = 411409.7740157063
$ws.Range("R4").Value = 7025032.483387623
# Row 10's record had no K/L/M/N/AC entries, so row 4 loses them.
$ws.Range("K4:N4").ClearContents()
$ws.Range("AC4").ClearContents()

# --- Row 5 (becomes what row 4 used to hold) ------------------------------
$ws.Range("A5").Value = 110060417
$ws.Range("Q5").Value = 411283.6994215957
$ws.Range("R5").Value = 7025288.901280839
$ws.Range("AC5").Value = "ringhack"

# --- Row 6 (becomes what row 5 used to hold) ------------------------------
$ws.Range("A6").Value = 110060420
$ws.Range("Q6").Value = 411272.1662275246
$ws.Range("R6").Value = 7025343.578228693

# --- Row 7 (becomes what row 6 used to hold) ------------------------------
$ws.Range("A7").Value = 110060415
$ws.Range("Q7").Value = 411290.2087285602
$ws.Range("R7").Value = 7025264.01279831

# --- Row 8 (becomes what row 7 used to hold) ------------------------------
$ws.Range("A8").Value = 110060425
$ws.Range("Q8").Value = 411500.4353926617
$ws.Range("R8").Value = 7025236.641728432
$ws.Range("AC8").Value = "ringhack färska"

# --- Row 9 (becomes what row 8 used to hold) ------------------------------
$ws.Range("A9").Value = 110060419
$ws.Range("B9").Value = 56395
$ws.Range("E9").Value = 100109
$ws.Range("F9").Value = "Tretåig hackspett"
$ws.Range("G9").Value = "Picoides tridactylus"
$ws.Range("H9").Value = "(Linnaeus, 1758)"
# Row 8's record had (empty) K/L/M/N entries, so row 9 gains them back.
# Copy the already-present-but-empty I9 cell into each so the cells exist
# (an empty-string .Value assignment would just delete the cell again).
$ws.Range("I9").Copy($ws.Range("K9"))
$ws.Range("I9").Copy($ws.Range("L9"))
$ws.Range("I9").Copy($ws.Range("M9"))
$ws.Range("I9").Copy($ws.Range("N9"))
$ws.Range("Q9").Value = 411283.4832894632
$ws.Range("R9").Value = 7025313.615521053
$ws.Range("AC9").Value = "ringhack"

# --- Row 10 (becomes what row 9 used to hold) -----------------------------
$ws.Range("A10").Value = 110060428
$ws.Range("B10").Value = 77506
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("Q10").Value = 411290.9159910702
$ws.Range("R10").Value = 7025370.912074795
